# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# with latest scraped values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.588.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0619"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.895.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.658.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.552.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "240.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0730"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.24%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.26%  "
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.111"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0502"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.461.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.925"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.94%  "
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.572"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "69.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.790"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.802.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("E49").Value = "  -5.99%  "
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.28%  "
